$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0) used between the currency symbol and amount,
# matching the formatting already used by the other price cells in column F.
$nbsp = [char]0x00A0

$ws.Range("F2").Value = "$" + $nbsp + "453,47"
$ws.Range("F4").Value = "$" + $nbsp + "1.851,15"
$ws.Range("F5").Value = "$" + $nbsp + "1.480,90"
$ws.Range("F11").Value = "$" + $nbsp + "444,20"
$ws.Range("F13").Value = "$" + $nbsp + "925,52"
$ws.Range("F17").Value = "$" + $nbsp + "3.054,45"
$ws.Range("F21").Value = "$" + $nbsp + "1.249,49"
$ws.Range("F22").Value = "Sin precio"
$ws.Range("F24").Value = "$" + $nbsp + "3.239,58"
$ws.Range("F25").Value = "$" + $nbsp + "1.295,77"
$ws.Range("F26").Value = "Sin precio"
$ws.Range("F27").Value = "$" + $nbsp + "416,43"
$ws.Range("F28").Value = "$" + $nbsp + "462,72"
$ws.Range("F29").Value = "$" + $nbsp + "462,72"
$ws.Range("F30").Value = "$" + $nbsp + "416,43"
$ws.Range("F31").Value = "$" + $nbsp + "832,97"
$ws.Range("F38").Value = "$" + $nbsp + "1.388,33"
$ws.Range("F39").Value = "$" + $nbsp + "499,74"
$ws.Range("F40").Value = "$" + $nbsp + "1.295,77"
$ws.Range("F41").Value = "$" + $nbsp + "990,32"
$ws.Range("F46").Value = "$" + $nbsp + "3.332,13"
$ws.Range("F48").Value = "$" + $nbsp + "638,58"
$ws.Range("F51").Value = "$" + $nbsp + "638,58"
$ws.Range("F52").Value = "$" + $nbsp + "879,24"
$ws.Range("F56").Value = "$" + $nbsp + "601,56"
$ws.Range("F57").Value = "$" + $nbsp + "601,56"
$ws.Range("F65").Value = "$" + $nbsp + "509,00"
$ws.Range("F66").Value = "$" + $nbsp + "509,00"
$ws.Range("F67").Value = "Sin precio"
$ws.Range("F72").Value = "$" + $nbsp + "1.129,16"
$ws.Range("F73").Value = "$" + $nbsp + "740,40"
$ws.Range("F75").Value = "$" + $nbsp + "2.406,52"
$ws.Range("F77").Value = "$" + $nbsp + "1.064,37"
$ws.Range("F78").Value = "$" + $nbsp + "1.064,37"
$ws.Range("F79").Value = "$" + $nbsp + "1.064,37"
$ws.Range("F80").Value = "$" + $nbsp + "1.064,37"
$ws.Range("F81").Value = "$" + $nbsp + "1.064,37"
$ws.Range("F82").Value = "$" + $nbsp + "2.961,88"
